# Append the new daily profit row (run date 2025-11-17) to the bottom of
# the data table on Sheet1, mirroring the existing layout:
#   column A -> date as literal text (e.g. "11/16/2025"), no special style
#   column B -> numeric profit value
#
# The new row goes immediately after the current last used row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.SpecialCells(11).Row   # xlCellTypeLastCell = 11
$newRow = $lastRow + 1

$dateCell = $ws.Cells.Item($newRow, 1)
$profitCell = $ws.Cells.Item($newRow, 2)

# Force column A to stay plain text instead of Excel auto-parsing the
# "MM/DD/YYYY" string into a date serial (which is how every prior row in
# this sheet is stored), then drop the temporary "@" format so the cell
# ends up with no explicit style, just like its neighbours.
$dateCell.NumberFormat = "@"
$dateCell.Value = "11/17/2025"
$dateCell.ClearFormats()

$profitCell.Value = 8827.700000000001
